$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells AD1:AF1, reusing the exact same formatting (bold, centered,
# top-aligned, thin box border) already used by the rest of row 1 (e.g. A1)
# by copying that cell's format/style onto the new header cells first.
$ws.Range("A1").Copy($ws.Range("AD1"))
$ws.Range("A1").Copy($ws.Range("AE1"))
$ws.Range("A1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the Wins/Losses/Ties data for rows 2-52
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 64
    $ws.Cells.Item($r, 31).Value = 97
    $ws.Cells.Item($r, 32).Value = 0
}
